$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Historico")

$row = 54

$ws.Cells.Item($row, 1).Value = "05/01/2026 08:52:03"
$ws.Cells.Item($row, 2).Value = "05/01 08:49"
$ws.Cells.Item($row, 3).Value = "g1 > Política"
$ws.Cells.Item($row, 4).Value = "Brasil vai condenar na ONU ataque à Venezuela, mas sabe que reunião não irá mudar situação no país vizinho"
$ws.Cells.Item($row, 5).Value = "https://g1.globo.com/politica/blog/valdo-cruz/post/2026/01/05/brasil-vai-condenar-na-onu-ataque-a-venezuela-mas-sabe-que-reuniao-nao-ira-mudar-situacao-no-pais-vizinho.ghtml"
$ws.Cells.Item($row, 6).Value = "lula"
$ws.Cells.Item($row, 7).Value = "s à Venezuela, mas sabe que a reunião não irá mudar a situação no país vizinho. O governo Lula decidiu fazer uso da palavra para reafirmar sua política tradicional em defesa da soberan"
